$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shifted rows 232-251 (D, L, M, N, O, P, S) ---
$ws.Range("D232").Value = 44461
$ws.Range("L232").Value = "1a amarillo"
$ws.Range("M232").Value = 300
$ws.Range("N232").Value = 5500
$ws.Range("O232").Value = 5500
$ws.Range("P232").Value = 5500
$ws.Range("S232").Value = 344

$ws.Range("D233").Value = 44461
$ws.Range("L233").Value = "2a amarillo"
$ws.Range("M233").Value = 300
$ws.Range("N233").Value = 4500
$ws.Range("O233").Value = 4500
$ws.Range("P233").Value = 4500
$ws.Range("S233").Value = 281

$ws.Range("D234").Value = 44357
$ws.Range("L234").Value = "1a amarillo"
$ws.Range("M234").Value = 400
$ws.Range("N234").Value = 8000
$ws.Range("O234").Value = 8000
$ws.Range("P234").Value = 8000
$ws.Range("S234").Value = 500

$ws.Range("D235").Value = 44357
$ws.Range("L235").Value = "2a amarillo"
$ws.Range("M235").Value = 300
$ws.Range("N235").Value = 6000
$ws.Range("O235").Value = 6000
$ws.Range("P235").Value = 6000
$ws.Range("S235").Value = 375

$ws.Range("D236").Value = 44203
$ws.Range("L236").Value = "1a amarillo"
$ws.Range("M236").Value = 300
$ws.Range("N236").Value = 17000
$ws.Range("O236").Value = 17000
$ws.Range("P236").Value = 17000
$ws.Range("S236").Value = 1062

$ws.Range("D237").Value = 44203
$ws.Range("L237").Value = "2a amarillo"
$ws.Range("M237").Value = 300
$ws.Range("N237").Value = 15000
$ws.Range("O237").Value = 15000
$ws.Range("P237").Value = 15000
$ws.Range("S237").Value = 938

$ws.Range("D238").Value = 44162
$ws.Range("L238").Value = "1a amarillo"
$ws.Range("M238").Value = 300
$ws.Range("N238").Value = 8000
$ws.Range("O238").Value = 8000
$ws.Range("P238").Value = 8000
$ws.Range("S238").Value = 500

$ws.Range("D239").Value = 44162
$ws.Range("L239").Value = "2a amarillo"
$ws.Range("M239").Value = 200
$ws.Range("N239").Value = 6500
$ws.Range("O239").Value = 6500
$ws.Range("P239").Value = 6500
$ws.Range("S239").Value = 406

$ws.Range("D240").Value = 44257
$ws.Range("L240").Value = "1a amarillo"
$ws.Range("M240").Value = 400
$ws.Range("N240").Value = 18000
$ws.Range("O240").Value = 18000
$ws.Range("P240").Value = 18000
$ws.Range("S240").Value = 1125

$ws.Range("D241").Value = 44257
$ws.Range("L241").Value = "2a amarillo"
$ws.Range("M241").Value = 300
$ws.Range("N241").Value = 16000
$ws.Range("O241").Value = 16000
$ws.Range("P241").Value = 16000
$ws.Range("S241").Value = 1000

$ws.Range("D242").Value = 44376
$ws.Range("L242").Value = "1a amarillo"
$ws.Range("M242").Value = 400
$ws.Range("N242").Value = 5500
$ws.Range("O242").Value = 5500
$ws.Range("P242").Value = 5500
$ws.Range("S242").Value = 344

$ws.Range("D243").Value = 44376
$ws.Range("L243").Value = "2a amarillo"
$ws.Range("M243").Value = 400
$ws.Range("N243").Value = 4500
$ws.Range("O243").Value = 4500
$ws.Range("P243").Value = 4500
$ws.Range("S243").Value = 281

$ws.Range("D244").Value = 44292
$ws.Range("L244").Value = "1a amarillo"
$ws.Range("M244").Value = 300
$ws.Range("N244").Value = 17000
$ws.Range("O244").Value = 17000
$ws.Range("P244").Value = 17000
$ws.Range("S244").Value = 1062

$ws.Range("D245").Value = 44292
$ws.Range("L245").Value = "2a amarillo"
$ws.Range("M245").Value = 300
$ws.Range("N245").Value = 15000
$ws.Range("O245").Value = 15000
$ws.Range("P245").Value = 15000
$ws.Range("S245").Value = 938

$ws.Range("D246").Value = 44358
$ws.Range("L246").Value = "1a amarillo"
$ws.Range("M246").Value = 400
$ws.Range("N246").Value = 8000
$ws.Range("O246").Value = 8000
$ws.Range("P246").Value = 8000
$ws.Range("S246").Value = 500

$ws.Range("D247").Value = 44358
$ws.Range("L247").Value = "2a amarillo"
$ws.Range("M247").Value = 400
$ws.Range("N247").Value = 6000
$ws.Range("O247").Value = 6000
$ws.Range("P247").Value = 6000
$ws.Range("S247").Value = 375

$ws.Range("D248").Value = 44211
$ws.Range("L248").Value = "1a plateado"
$ws.Range("M248").Value = 600
$ws.Range("N248").Value = 20000
$ws.Range("O248").Value = 21000
$ws.Range("P248").Value = 20500
$ws.Range("S248").Value = 1281

$ws.Range("D249").Value = 44211
$ws.Range("L249").Value = "2a plateado"
$ws.Range("M249").Value = 300
$ws.Range("N249").Value = 18000
$ws.Range("O249").Value = 18000
$ws.Range("P249").Value = 18000
$ws.Range("S249").Value = 1125

$ws.Range("D250").Value = 44425
$ws.Range("L250").Value = "1a amarillo"
$ws.Range("M250").Value = 400
$ws.Range("N250").Value = 5000
$ws.Range("O250").Value = 5000
$ws.Range("P250").Value = 5000
$ws.Range("S250").Value = 312

$ws.Range("D251").Value = 44425
$ws.Range("L251").Value = "2a amarillo"
$ws.Range("M251").Value = 300
$ws.Range("N251").Value = 4000
$ws.Range("O251").Value = 4000
$ws.Range("P251").Value = 4000
$ws.Range("S251").Value = 250

# --- New rows 252 and 253 (full rows, same fixed columns as the block above) ---
$ws.Range("A252").Value = 11
$ws.Range("B252").Value = "Vega Monumental Concepción"
$ws.Range("C252").Value = "Bíobío"
$ws.Range("D252").Value = 44323
$ws.Range("E252").Value = 8
$ws.Range("F252").Value = "Fruta"
$ws.Range("G252").Value = 100102
$ws.Range("H252").Value = "Cítricos"
$ws.Range("I252").Value = 100102003
$ws.Range("J252").Value = "Limón"
$ws.Range("K252").Value = "Sin especificar"
$ws.Range("L252").Value = "1a amarillo"
$ws.Range("M252").Value = 300
$ws.Range("N252").Value = 14000
$ws.Range("O252").Value = 14000
$ws.Range("P252").Value = 14000
$ws.Range("Q252").Value = "$/malla 16 kilos"
$ws.Range("R252").Value = "Provincia de Melipilla"
$ws.Range("S252").Value = 875
$ws.Range("T252").Value = 16
$ws.Range("D252").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A253").Value = 11
$ws.Range("B253").Value = "Vega Monumental Concepción"
$ws.Range("C253").Value = "Bíobío"
$ws.Range("D253").Value = 44323
$ws.Range("E253").Value = 8
$ws.Range("F253").Value = "Fruta"
$ws.Range("G253").Value = 100102
$ws.Range("H253").Value = "Cítricos"
$ws.Range("I253").Value = 100102003
$ws.Range("J253").Value = "Limón"
$ws.Range("K253").Value = "Sin especificar"
$ws.Range("L253").Value = "2a amarillo"
$ws.Range("M253").Value = 300
$ws.Range("N253").Value = 12000
$ws.Range("O253").Value = 12000
$ws.Range("P253").Value = 12000
$ws.Range("Q253").Value = "$/malla 16 kilos"
$ws.Range("R253").Value = "Provincia de Melipilla"
$ws.Range("S253").Value = 750
$ws.Range("T253").Value = 16
$ws.Range("D253").NumberFormat = "YYYY-MM-DD HH:MM:SS"

